{"js": "// Replace the 100 arithmetic-expression cells in the single table with\n// their updated expressions. The mapping below is the ordered list of\n// (old, new) text pairs exactly as they appear, left-to-right / top-to-\n// bottom, in the document's one table (20 rows x 5 columns). Several old\n// values repeat (e.g. \"88-26=\" and \"72-40=\" each appear twice), so we\n// replace using the saved per-cell \"new\" values directly rather than a\n// text search, which keeps the mapping unambiguous and position-exact.\nconst newValues = [\n  [\"19+69=\", \"20+4=\", \"79+6=\", \"7+42=\", \"66+12=\"],\n  [\"13+51=\", \"97-65=\", \"73+21=\", \"70-1=\", \"29+24=\"],\n  [\"41+35=\", \"21+78=\", \"86-23=\", \"33+20=\", \"63-62=\"],\n  [\"73-17=\", \"37-20=\", \"19+45=\", \"22+6=\", \"63-45=\"],\n  [\"84-65=\", \"41-8=\", \"71-2=\", \"34+23=\", \"68-24=\"],\n  [\"84+10=\", \"9+35=\", \"5+20=\", \"28-10=\", \"97-29=\"],\n  [\"67+10=\", \"99-44=\", \"14+38=\", \"85-50=\", \"73-15=\"],\n  [\"20-16=\", \"7+63=\", \"58+32=\", \"51-39=\", \"78+7=\"],\n  [\"0+99=\", \"70-30=\", \"42+40=\", \"40-22=\", \"57-20=\"],\n  [\"80-13=\", \"34-7=\", \"35+8=\", \"81-15=\", \"50-4=\"],\n  [\"56+40=\", \"21+32=\", \"33+45=\", \"1+30=\", \"8+60=\"],\n  [\"68+22=\", \"57-7=\", \"19+0=\", \"93-35=\", \"93-4=\"],\n  [\"71-61=\", \"33-31=\", \"62-24=\", \"76+10=\", \"76-60=\"],\n  [\"75+1=\", \"10+75=\", \"5+66=\", \"96-74=\", \"85-68=\"],\n  [\"52-30=\", \"38+26=\", \"92-55=\", \"93-21=\", \"24+20=\"],\n  [\"45-29=\", \"49+14=\", \"91+8=\", \"16-6=\", \"61+32=\"],\n  [\"22+77=\", \"14+83=\", \"20+45=\", \"54+31=\", \"18+10=\"],\n  [\"36+22=\", \"28+15=\", \"21-8=\", \"50-23=\", \"86+0=\"],\n  [\"59+2=\", \"5+9=\", \"42+33=\", \"99-49=\", \"98-55=\"],\n  [\"9+60=\", \"60-14=\", \"47+31=\", \"78-26=\", \"87-34=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// Build the replacement grid from the current cell contents so formatting-\n// only rows/columns (if any ever existed) are left untouched, and only\n// assign where we have a mapped replacement.\nconst updated = table.values.map((row, r) =>\n  row.map((cell, c) => {\n    if (newValues[r] && newValues[r][c] !== undefined) {\n      return newValues[r][c];\n    }\n    return cell;\n  })\n);\n\ntable.values = updated;\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-expression cells in the single table with\n# their updated expressions. $newValues is the ordered grid of\n# replacement text, row-major (20 rows x 5 columns), matching the\n# table's cell order top-to-bottom / left-to-right.\n$newValues = @(\n    @(\"19+69=\", \"20+4=\", \"79+6=\", \"7+42=\", \"66+12=\"),\n    @(\"13+51=\", \"97-65=\", \"73+21=\", \"70-1=\", \"29+24=\"),\n    @(\"41+35=\", \"21+78=\", \"86-23=\", \"33+20=\", \"63-62=\"),\n    @(\"73-17=\", \"37-20=\", \"19+45=\", \"22+6=\", \"63-45=\"),\n    @(\"84-65=\", \"41-8=\", \"71-2=\", \"34+23=\", \"68-24=\"),\n    @(\"84+10=\", \"9+35=\", \"5+20=\", \"28-10=\", \"97-29=\"),\n    @(\"67+10=\", \"99-44=\", \"14+38=\", \"85-50=\", \"73-15=\"),\n    @(\"20-16=\", \"7+63=\", \"58+32=\", \"51-39=\", \"78+7=\"),\n    @(\"0+99=\", \"70-30=\", \"42+40=\", \"40-22=\", \"57-20=\"),\n    @(\"80-13=\", \"34-7=\", \"35+8=\", \"81-15=\", \"50-4=\"),\n    @(\"56+40=\", \"21+32=\", \"33+45=\", \"1+30=\", \"8+60=\"),\n    @(\"68+22=\", \"57-7=\", \"19+0=\", \"93-35=\", \"93-4=\"),\n    @(\"71-61=\", \"33-31=\", \"62-24=\", \"76+10=\", \"76-60=\"),\n    @(\"75+1=\", \"10+75=\", \"5+66=\", \"96-74=\", \"85-68=\"),\n    @(\"52-30=\", \"38+26=\", \"92-55=\", \"93-21=\", \"24+20=\"),\n    @(\"45-29=\", \"49+14=\", \"91+8=\", \"16-6=\", \"61+32=\"),\n    @(\"22+77=\", \"14+83=\", \"20+45=\", \"54+31=\", \"18+10=\"),\n    @(\"36+22=\", \"28+15=\", \"21-8=\", \"50-23=\", \"86+0=\"),\n    @(\"59+2=\", \"5+9=\", \"42+33=\", \"99-49=\", \"98-55=\"),\n    @(\"9+60=\", \"60-14=\", \"47+31=\", \"78-26=\", \"87-34=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
